# Update workbook data (gh-pages output regenerated) across all four sheets:
# 展览 (Exhibition), 演出 (Performance), 本地生活 (Local Life), 全部类型 (All Types)
#
# Column F = "想去人数" (want-to-go count), Column G = "最低票价" (lowest ticket price)
# For "广州·COMICUP 2024SP" the event is now sold out, so its F becomes a new
# date-serial-looking counter value and G switches from a numeric price to the
# text "已售罄" ("Sold Out").

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 37712
$ws.Range("G2").Value = "已售罄"
$ws.Range("F4").Value = 638
$ws.Range("F5").Value = 777
$ws.Range("F8").Value = 469
$ws.Range("F9").Value = 851
$ws.Range("F11").Value = 721
$ws.Range("F12").Value = 560
$ws.Range("F13").Value = 55
$ws.Range("F14").Value = 37
$ws.Range("F15").Value = 25
$ws.Range("F16").Value = 658
$ws.Range("F17").Value = 181
$ws.Range("F19").Value = 444
$ws.Range("F20").Value = 1174
$ws.Range("F22").Value = 839
$ws.Range("F23").Value = 2544
$ws.Range("F24").Value = 1023
$ws.Range("F25").Value = 567
$ws.Range("F27").Value = 1165
$ws.Range("F29").Value = 787
$ws.Range("F30").Value = 65
$ws.Range("F31").Value = 1166

# ---- Sheet 2: 演出 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 411
$ws.Range("F4").Value = 333
$ws.Range("F10").Value = 12
$ws.Range("F11").Value = 10

# ---- Sheet 3: 本地生活 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 638

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 638
$ws.Range("F3").Value = 37712
$ws.Range("G3").Value = "已售罄"
$ws.Range("F5").Value = 638
$ws.Range("F6").Value = 777
$ws.Range("F10").Value = 469
$ws.Range("F11").Value = 411
$ws.Range("F12").Value = 333
$ws.Range("F15").Value = 851
$ws.Range("F17").Value = 721
$ws.Range("F18").Value = 560
$ws.Range("F19").Value = 55
$ws.Range("F21").Value = 37
$ws.Range("F24").Value = 12
$ws.Range("F25").Value = 25
$ws.Range("F26").Value = 10
$ws.Range("F27").Value = 658
$ws.Range("F28").Value = 181
$ws.Range("F30").Value = 444
$ws.Range("F31").Value = 1174
$ws.Range("F33").Value = 839
$ws.Range("F34").Value = 2544
$ws.Range("F35").Value = 1023
$ws.Range("F36").Value = 567
$ws.Range("F38").Value = 1165
$ws.Range("F41").Value = 787
$ws.Range("F42").Value = 65
$ws.Range("F43").Value = 1166

$wb.Save()
